$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28 (pushes old rows 28-42 down to 29-43),
# carrying down formatting from the row below as Excel normally does.
$ws.Rows.Item(28).Insert()

# Fill in the new row 28 with the new data record.
$ws.Range("A28").Value = 11
$ws.Range("B28").Value = "Vega Monumental Concepción"
$ws.Range("C28").Value = "Bíobío"
$ws.Range("D28").Value = 44460
$ws.Range("D28").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E28").Value = 8
$ws.Range("F28").Value = 100112021
$ws.Range("G28").Value = "Ají"
$ws.Range("H28").Value = "Inferno"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 50
$ws.Range("K28").Value = 35000
$ws.Range("L28").Value = 36000
$ws.Range("M28").Value = 35400
$ws.Range("N28").Value = "$/caja 12 kilos"
$ws.Range("O28").Value = "Región de Arica y Parinacota"
$ws.Range("P28").Value = 2950
$ws.Range("Q28").Value = 12
$ws.Range("R28").Value = "Hortaliza"
